$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.5
$ws.Range("I5").Value = 3.3
$ws.Range("J5").Value = 1.1
$ws.Range("L5").Value = 1.58
$ws.Range("M5").Value = 2.2
$ws.Range("N5").Value = 3.1
$ws.Range("O5").Value = 1.33
$ws.Range("T5").Value = 6
$ws.Range("U5").Value = 10
$ws.Range("W5").Value = 26
$ws.Range("AF5").Value = 13
$ws.Range("AH5").Value = 34
$ws.Range("AI5").Value = 34

# Row 6
$ws.Range("G6").Value = 3.9
$ws.Range("I6").Value = 2.15
$ws.Range("T6").Value = 7.5
$ws.Range("U6").Value = 17
$ws.Range("V6").Value = 15
$ws.Range("X6").Value = 41
$ws.Range("AF6").Value = 8.5
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 19

# Row 10
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 4.5
$ws.Range("N10").Value = 1.67
$ws.Range("O10").Value = 2.15
$ws.Range("U10").Value = 9
$ws.Range("W10").Value = 13
$ws.Range("Y10").Value = 21
$ws.Range("AA10").Value = 7.5
$ws.Range("AF10").Value = 26
$ws.Range("AH10").Value = 51

# Row 17
$ws.Range("G17").Value = 2.8
$ws.Range("H17").Value = 2.8
$ws.Range("I17").Value = 2.8
$ws.Range("T17").Value = 7
$ws.Range("U17").Value = 12
$ws.Range("W17").Value = 29
$ws.Range("X17").Value = 26
$ws.Range("AE17").Value = 7
$ws.Range("AF17").Value = 12
$ws.Range("AG17").Value = 11
$ws.Range("AH17").Value = 29
$ws.Range("AI17").Value = 26

# Row 18
$ws.Range("G18").Value = 1.83
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 5.25
$ws.Range("L18").Value = 1.53
$ws.Range("M18").Value = 2.38
$ws.Range("N18").Value = 2.7
$ws.Range("O18").Value = 1.44
$ws.Range("P18").Value = 1.57
$ws.Range("Q18").Value = 2.25
$ws.Range("V18").Value = 8
$ws.Range("Z18").Value = 6
$ws.Range("AA18").Value = 6
$ws.Range("AC18").Value = 81

# Row 20
$ws.Range("L20").Value = 1.34
$ws.Range("M20").Value = 2.75
$ws.Range("N20").Value = 2
$ws.Range("P20").Value = 1.44
$ws.Range("Q20").Value = 2.42
$ws.Range("R20").Value = 2.1
$ws.Range("S20").Value = 1.57
$ws.Range("T20").Value = 5.5
$ws.Range("U20").Value = 6.2
$ws.Range("Y20").Value = 37
$ws.Range("Z20").Value = 8.75
$ws.Range("AC20").Value = 150
$ws.Range("AG20").Value = 19.5
$ws.Range("AI20").Value = 75

# Row 22
$ws.Range("J22").Value = 1.07
$ws.Range("K22").Value = 9

# Row 25
$ws.Range("G25").Value = 1.73
$ws.Range("I25").Value = 4.5
$ws.Range("L25").Value = 1.33
$ws.Range("M25").Value = 3.25
$ws.Range("N25").Value = 2.08
$ws.Range("O25").Value = 1.73
$ws.Range("T25").Value = 6
$ws.Range("U25").Value = 7.5
$ws.Range("V25").Value = 8.5
$ws.Range("W25").Value = 13
$ws.Range("X25").Value = 15
$ws.Range("AE25").Value = 11
$ws.Range("AF25").Value = 23
$ws.Range("AG25").Value = 15
$ws.Range("AH25").Value = 51
$ws.Range("AI25").Value = 41

# Row 60
$ws.Range("L60").Value = 1.18
$ws.Range("M60").Value = 4.5
$ws.Range("N60").Value = 1.62
$ws.Range("O60").Value = 2.25
$ws.Range("R60").Value = 1.5

# Row 61
$ws.Range("G61").Value = 1.7
$ws.Range("H61").Value = 3.65
$ws.Range("I61").Value = 4.4
$ws.Range("L61").Value = 1.18
$ws.Range("M61").Value = 4.5
$ws.Range("R61").Value = 1.53
$ws.Range("S61").Value = 2.18
$ws.Range("U61").Value = 9.75
$ws.Range("W61").Value = 14.5
$ws.Range("X61").Value = 12
$ws.Range("Y61").Value = 19
$ws.Range("Z61").Value = 13.5
$ws.Range("AA61").Value = 7.4
$ws.Range("AB61").Value = 12.5
$ws.Range("AC61").Value = 45
$ws.Range("AD61").Value = 250
$ws.Range("AE61").Value = 16
$ws.Range("AF61").Value = 30
$ws.Range("AG61").Value = 14
$ws.Range("AH61").Value = 80
$ws.Range("AI61").Value = 37
$ws.Range("AJ61").Value = 35

# Row 65
$ws.Range("G65").Value = 1.83
$ws.Range("H65").Value = 3
$ws.Range("I65").Value = 4.33
$ws.Range("J65").Value = 1.07
$ws.Range("L65").Value = 1.47
$ws.Range("P65").Value = 1.57
$ws.Range("Q65").Value = 2.25
$ws.Range("R65").Value = 2.2
$ws.Range("S65").Value = 1.62
$ws.Range("U65").Value = 7.5
$ws.Range("W65").Value = 15
$ws.Range("AE65").Value = 9.5
$ws.Range("AF65").Value = 21

# Row 86
$ws.Range("G86").Value = 1.65
$ws.Range("H86").Value = 3.55
$ws.Range("I86").Value = 5
$ws.Range("L86").Value = 1.26
$ws.Range("M86").Value = 3.15
$ws.Range("N86").Value = 1.78
$ws.Range("O86").Value = 1.83
$ws.Range("P86").Value = 1.39
$ws.Range("Q86").Value = 2.57
$ws.Range("R86").Value = 1.75
$ws.Range("S86").Value = 1.87
$ws.Range("T86").Value = 7.1
$ws.Range("U86").Value = 8
$ws.Range("V86").Value = 7.9
$ws.Range("W86").Value = 13
$ws.Range("X86").Value = 12.5
$ws.Range("Y86").Value = 24
$ws.Range("Z86").Value = 10
$ws.Range("AA86").Value = 6.9
$ws.Range("AB86").Value = 15.5
$ws.Range("AC86").Value = 70
$ws.Range("AD86").Value = 600
$ws.Range("AE86").Value = 13
$ws.Range("AF86").Value = 30
$ws.Range("AG86").Value = 16
$ws.Range("AH86").Value = 100
$ws.Range("AI86").Value = 55
$ws.Range("AJ86").Value = 55
